# Data.xlsx ("Ticketonline/TestData/Data.xlsx") — commit "commited code on 20 oct 2019"
#
# Semantic changes in the sheet itself:
#   1. The password stored in B1 ("Test@123") is replaced by "Arkadmin@1".
#   2. The sheet's last-saved selection moves from B8 to B2.
#
# (The workbookView window-chrome numbers, the calcPr/calcId value and the
#  <oleSize> element are pure Excel-UI/session bookkeeping baked in by the
#  authoring client when the file was last closed — they are not reachable
#  through the Excel object model and are left to the host application.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the credential stored in B1.
$ws.Range("B1").Value = "Arkadmin@1"

# 2. Move the active selection from B8 to B2.
$ws.Range("B2").Select()
